$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("singlish_NUS_SMS")

$ws.Range("C6").Value = "`n`nA: There is duty free sales on the plane, you buy from here lah`nB: It's fine, no thank you`nA: The prices are better here`nB: Are you sure?`nA: Ya, confirm plus chop"

$ws.Range("C10").Value = "`n`nA: Hey, what's up?`nB: My baby sneezed several times this morning.`nA: Oh no, is he sick?`nB: I don't know, but my wife insisted I should take him to the doctor leh.`nA: I think he's just fine, tell your wife don't need to be so anxious"

$ws.Activate()
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 14
